# Atualização automática de preços de eletricidade
# Daily automatic update: row 2 (the only data row) moves forward one day
# and all the hourly/slot prices are refreshed with the new day's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45949

$ws.Range("B2").Value = 86.54000000000001
$ws.Range("C2").Value = 71.36
$ws.Range("D2").Value = 56.27
$ws.Range("E2").Value = 34.13
$ws.Range("F2").Value = 30.47
$ws.Range("G2").Value = 29.58
$ws.Range("H2").Value = 29.75
$ws.Range("I2").Value = 30.62
$ws.Range("J2").Value = 31.35
$ws.Range("K2").Value = 24.8
$ws.Range("L2").Value = 5.88
$ws.Range("M2").Value = 1.67
$ws.Range("N2").Value = 0.01
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.17
$ws.Range("S2").Value = 3.46
$ws.Range("T2").Value = 36.2
$ws.Range("U2").Value = 66.65000000000001
$ws.Range("V2").Value = 72.14
$ws.Range("W2").Value = 70.86
$ws.Range("X2").Value = 57.38
$ws.Range("Y2").Value = 47.41
$ws.Range("Z2").Value = 32.78

$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 62.08
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 78.95
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 71.5
$ws.Range("AG2").Value = "4h-17h"
